$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D ("Celebrante"), shifting
# Celebrante -> E and Ministros_Eucaristia -> F. The new column becomes the
# "Credência" column.
$ws.Columns("D").Insert()

# Header for the new column.
$ws.Range("D1").Value = "Credência"

# Cyclic roster of names for the new "Credência" column (rows 2-28).
$names = @("Mario", "Jorge", "José", "Carlos", "Simone")
for ($r = 2; $r -le 28; $r++) {
    $idx = ($r - 2) % 5
    $ws.Cells.Item($r, 4).Value = $names[$idx]
}

# Match the column width used in the edited workbook as closely as this
# engine's width model allows.
$ws.Columns("D").ColumnWidth = 8.498697916666666

# Update the selection to match the edited workbook.
$ws.Range("D28").Select() | Out-Null
